$wb = $excel.ActiveWorkbook

# --- Currencies sheet: selection only changes (no data edits) ---
$wsCurrencies = $wb.Worksheets.Item("Currencies")
$wsCurrencies.Activate()
$wsCurrencies.Range("F8").Select()

# --- Payee sheet: add a new "JSON row" helper column (J) ---
$wsPayee = $wb.Worksheets.Item("Payee")
$wsPayee.Activate()

# First data row (row 4) gets its own (non-shared) formula, matching the
# existing pattern used by columns E/F/H on this sheet.
$wsPayee.Range("J4").Formula = '="{" & CHAR(34) & $C4 & CHAR(34) & "," & CHAR(34) &  $D4 & CHAR(34) & "},"'

# Remaining rows (5-24) share one formula, auto-filled from row 5's formula.
$wsPayee.Range("J5:J24").Formula = '="{" & CHAR(34) & $C5 & CHAR(34) & "," & CHAR(34) &  $D5 & CHAR(34) & "},"'

# Touch page setup so this sheet gets an explicit <pageSetup> entry, as on
# the Currencies sheet.
$psPayee = $wsPayee.PageSetup
$psPayee.PaperSize = 9
$psPayee.Orientation = 1

$wsPayee.Range("H4:H24").Select()

# --- Recurring Charge sheet: no longer the active sheet; selection widens ---
# (Select()-ing a range on a non-active sheet activates it, so do this
# before activating Account, which must be the final active sheet.)
$wsRecurring = $wb.Worksheets.Item("Recurring Charge")
$wsRecurring.Activate()
$wsRecurring.Range("M4:M16").Select()

# --- Account sheet: becomes the active sheet; selection widens to K4:K6 ---
$wsAccount = $wb.Worksheets.Item("Account")
$wsAccount.Activate()
$wsAccount.Range("K4:K6").Select()
